# Feuille de temps - update for "4-sept" sheet and active-tab selection.
$wb = $excel.ActiveWorkbook

# --- Update the data rows on the "4-sept" sheet (A: Nom, C: Heure début, E: Tâche) ---
$ws4 = $wb.Worksheets.Item("4-sept")

# Row 3: Hugo / Content Matrix -> Gabriel / makette
$ws4.Range("A3").Value = "Gabriel"
$ws4.Range("E3").Value = "makette"

# Row 4: Gabriel -> Louis (time/task unchanged)
$ws4.Range("A4").Value = "Louis"

# Row 5: Louis / 0.35416666666666669 / makette -> Sophie / 0.36458333333333331 / Planification phase 1
$ws4.Range("A5").Value = "Sophie"
$ws4.Range("C5").Value = 0.36458333333333331
$ws4.Range("E5").Value = "Planification phase 1"

# Row 6: Sophie / 0.36458333333333331 / Planification phase 1 -> Hugo / 0.35416666666666669 / Plan BD
$ws4.Range("A6").Value = "Hugo"
$ws4.Range("C6").Value = 0.35416666666666669
$ws4.Range("E6").Value = "Plan BD"

# --- Switch the active/selected tab from "3-sept" to "4-sept" ---
$ws4.Activate()
$ws4.Range("G7").Select()
